# CGC : Tax Logic Fixes
# Updates the Short Term / Long Term capital-gain tax tables with corrected
# sample figures and makes the tax formulas defensive against negative
# gains (IF(...>=0 / >threshold ..., 0)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
[void]$ws.Activate()

# ---------------------------------------------------------------
# Short Term section (rows 3-5)
# ---------------------------------------------------------------

# Row 3: "Before 23rd July, 2024"
$ws.Range("B3").Value = 80000
$ws.Range("C3").Value = 75000
$ws.Range("D3").Formula = "=IF(B3-C3>=0,ROUND((B3-C3)*15%,0),0)"

# Row 4: "After 23rd July,2024"
$ws.Range("B4").Value = 80000
$ws.Range("C4").Value = 90000
$ws.Range("D4").Formula = "=IF(B4-C4>=0,ROUND((B4-C4)*20%,0),0)"

# ---------------------------------------------------------------
# Long Term section (rows 8-10)
# ---------------------------------------------------------------

# Row 9: "After 23rd July,2024"
$ws.Range("B9").Value = 300000
$ws.Range("C9").Value = 500000
$ws.Range("D9").Formula = "=IF(B9-C9>125000,ROUND(((B9-C9)-125000)*12.5%,0),0)"

# ---------------------------------------------------------------
# View state: scroll down one row and move the selection to the
# Long Term "Tax" header cell (E7:F7).
# ---------------------------------------------------------------

$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E7:F7").Select()
